$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price (D) and volume-change (E) columns per latest scrape.
# Rows 45/46 additionally swap the Coin name (B) and Link (C) — EnergySwap and
# Decentraland traded ranking positions.

$ws.Range("D2").Value = '27.843.00'
$ws.Range("E2").Value = '  -0.63%  '

$ws.Range("D3").Value = '1.906.75'
$ws.Range("E3").Value = '  -0.03%  '

$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").Value = '''313.41'
$ws.Range("E5").Value = '  -0.79%  '

$ws.Range("D6").Value = '''1.003'

$ws.Range("D7").Value = '''0.5030'
$ws.Range("E7").Value = '  +4.28%  '

$ws.Range("D8").Value = '''0.3809'
$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("D9").Value = '''0.07273'
$ws.Range("E9").Value = '  -1.07%  '

$ws.Range("D10").Value = '''0.9060'
$ws.Range("E10").Value = '  -2.86%  '

$ws.Range("D11").Value = '''20.85'
$ws.Range("E11").Value = '  +0.13%  '

$ws.Range("D12").Value = '''0.07661'
$ws.Range("E12").Value = '  -1.16%  '

$ws.Range("D13").Value = '1.909.95'
$ws.Range("E13").Value = '  -0.16%  '

$ws.Range("D14").Value = '''5.484'
$ws.Range("E14").Value = '  -0.36%  '

$ws.Range("D15").Value = '''91.67'
$ws.Range("E15").Value = '  -0.07%  '

$ws.Range("D16").Value = '''1.003'
$ws.Range("E16").Value = '  -0.35%  '

$ws.Range("D17").Value = '''0.000008719'
$ws.Range("E17").Value = '  -1.19%  '

$ws.Range("E18").Value = '  -0.24%  '

$ws.Range("D19").Value = '27.869.18'
$ws.Range("E19").Value = '  -0.67%  '

$ws.Range("D20").Value = '''14.59'
$ws.Range("E20").Value = '  -1.45%  '

$ws.Range("E21").Value = '  -0.28%  '

$ws.Range("D22").Value = '''10.81'
$ws.Range("E22").Value = '  -0.75%  '

$ws.Range("D23").Value = '''6.571'
$ws.Range("E23").Value = '  -1.01%  '

$ws.Range("D24").Value = '''153.67'
$ws.Range("E24").Value = '  -1.26%  '

$ws.Range("D25").Value = '''1.875'
$ws.Range("E25").Value = '  -2.32%  '

$ws.Range("D26").Value = '''2.216'
$ws.Range("E26").Value = '  +4.12%  '

$ws.Range("D27").Value = '''18.37'
$ws.Range("E27").Value = '  -0.73%  '

$ws.Range("D28").Value = '''115.27'
$ws.Range("E28").Value = '  -1.42%  '

$ws.Range("D29").Value = '''4.900'
$ws.Range("E29").Value = '  -1.22%  '

$ws.Range("E30").Value = '  +0.72%  '

$ws.Range("D31").Value = '''3.209'
$ws.Range("E31").Value = '  -2.91%  '

$ws.Range("D32").Value = '''1.220'
$ws.Range("E32").Value = '  -2.69%  '

$ws.Range("D33").Value = '''4.666'
$ws.Range("E33").Value = '  -0.18%  '

$ws.Range("D34").Value = '''0.7626'
$ws.Range("E34").Value = '  -1.65%  '

$ws.Range("D35").Value = '''0.02064'
$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("D36").Value = '''2.503'
$ws.Range("E36").Value = '  -4.69%  '

$ws.Range("D37").Value = '''1.093'
$ws.Range("E37").Value = '  -1.65%  '

$ws.Range("D38").Value = '''0.5524'
$ws.Range("E38").Value = '  +0.55%  '

$ws.Range("D39").Value = '''3.013'
$ws.Range("E39").Value = '  +0.65%  '

$ws.Range("D40").Value = '''0.05245'
$ws.Range("E40").Value = '  -1.18%  '

$ws.Range("D41").Value = '''6.860'
$ws.Range("E41").Value = '  -2.45%  '

$ws.Range("E42").Value = '  -0.71%  '

$ws.Range("D43").Value = '''0.1506'
$ws.Range("E43").Value = '  -1.84%  '

$ws.Range("D44").Value = '''111.10'
$ws.Range("E44").Value = '  +3.14%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '''0.4799'
$ws.Range("E45").Value = '  -0.69%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''10.52'
$ws.Range("E46").Value = '  -1.96%  '

$ws.Range("D47").Value = '''1.003'
$ws.Range("E47").Value = '  -0.13%  '

$ws.Range("D48").Value = '''1.622'
$ws.Range("E48").Value = '  -1.52%  '

$ws.Range("D49").Value = '''67.24'
$ws.Range("E49").Value = '  -1.18%  '

$ws.Range("D50").Value = '''0.06060'
$ws.Range("E50").Value = '  -0.13%  '

$ws.Range("D51").Value = '''0.9014'
$ws.Range("E51").Value = '  +0.18%  '
